$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("LEFT RIGHT MID Functions")

$ws.Range("G3").Formula = "=RIGHT(A3,LEN(A3)-6)"
$ws.Range("H3").Formula = "=LEN(A3)"

$ws.Range("G4:G26").Formula = "=RIGHT(A4,LEN(A4)-6)"
$ws.Range("H4:H26").Formula = "=LEN(A4)"
$ws.Range("E5:E26").Formula = "=LEFT(A5,3)"
$ws.Range("F5:F26").Formula = "=MID(A5,4,3)"

$ws.Range("B16").Select()
